$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21; this shifts the existing rows 21-25
# (Bíobío / Vega Monumental Concepción weekly Alcachofa entries) down to 22-26,
# preserving all of their original data.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the latest weekly price entry.
$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value = "Bíobío"
$ws.Cells.Item(21, 4).Value = 44463
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112013
$ws.Cells.Item(21, 7).Value = "Alcachofa"
$ws.Cells.Item(21, 8).Value = "Argentina(o)"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 9000
$ws.Cells.Item(21, 12).Value = 10000
$ws.Cells.Item(21, 13).Value = 9500
$ws.Cells.Item(21, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(21, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(21, 16).Value = 190
$ws.Cells.Item(21, 17).Value = 50
$ws.Cells.Item(21, 18).Value = "Hortaliza"
